$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 797.2857
$ws.Range("I19").Value = 916.3
$ws.Range("J19").Value = 499.75
$ws.Range("K19").Value = 916.3
$ws.Range("L19").Value = 499.75
$ws.Range("M19").Value = -741.3
$ws.Range("N19").Value = -849.75
$ws.Range("H76").Value = 12682.091
$ws.Range("J76").Value = 5000
$ws.Range("L76").Value = 5000
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 12682.091
$ws.Range("J79").Value = 5000
$ws.Range("L79").Value = 5000
$ws.Range("N79").Value = -7184
$ws.Range("H92").Value = 3167.8333
$ws.Range("I92").Value = 3200.4
$ws.Range("K92").Value = 3200.4
$ws.Range("M92").Value = -1952.4
$ws.Range("H113").Value = 8213.429
$ws.Range("J113").Value = 9497.5
$ws.Range("L113").Value = 9497.5
$ws.Range("N113").Value = -16005.5
$ws.Range("H132").Value = 1347.3135
$ws.Range("I132").Value = 1340.4546
$ws.Range("J132").Value = 1800
$ws.Range("K132").Value = 4021.3638
$ws.Range("L132").Value = 5400
$ws.Range("M132").Value = -1491.3638
$ws.Range("N132").Value = -10460
$ws.Range("H134").Value = 104166.5
$ws.Range("J134").Value = 104166.5
$ws.Range("L134").Value = 104166.5
$ws.Range("N134").Value = -114306.5
$ws.Range("H135").Value = 13889599
$ws.Range("I135").Value = 674.0625
$ws.Range("J135").Value = 125001000
$ws.Range("K135").Value = 6066.5625
$ws.Range("L135").Value = 1125009000
$ws.Range("M135").Value = -3531.5625
$ws.Range("N135").Value = -1125014070
$ws.Range("H137").Value = 4423.397
$ws.Range("I137").Value = 2017.4584
$ws.Range("J137").Value = 12122.4
$ws.Range("K137").Value = 6052.3752
$ws.Range("L137").Value = 36367.2
$ws.Range("M137").Value = -3502.3752
$ws.Range("N137").Value = -41467.2
$ws.Range("H138").Value = 23813520
$ws.Range("J138").Value = 66675948
$ws.Range("L138").Value = 200027844
$ws.Range("N138").Value = -200038124

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16955932
$ws.Range("I32").Value = 19235644
$ws.Range("J32").Value = 20933.857
$ws.Range("K32").Value = 19235644
$ws.Range("L32").Value = 20933.857
$ws.Range("M32").Value = -19235357
$ws.Range("N32").Value = -21507.857
$ws.Range("H45").Value = 2182.3704
$ws.Range("I45").Value = 2325.2104
$ws.Range("K45").Value = 2325.2104
$ws.Range("M45").Value = -1948.2104
$ws.Range("H74").Value = 41714580
$ws.Range("I74").Value = 62570868
$ws.Range("K74").Value = 62570868
$ws.Range("M74").Value = -62569994
$ws.Range("H77").Value = 41714580
$ws.Range("I77").Value = 62570868
$ws.Range("K77").Value = 312854340
$ws.Range("M77").Value = -312849972
$ws.Range("H135").Value = 190142
$ws.Range("J135").Value = 190142
$ws.Range("L135").Value = 190142
$ws.Range("N135").Value = -200282

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1536.75
$ws.Range("I20").Value = 1480.1
$ws.Range("J20").Value = 1820
$ws.Range("K20").Value = 1480.1
$ws.Range("L20").Value = 1820
$ws.Range("M20").Value = -1233.1
$ws.Range("N20").Value = -2314
$ws.Range("H59").Value = 119989.664
$ws.Range("J59").Value = 119989.664
$ws.Range("L59").Value = 119989.664
$ws.Range("N59").Value = -121683.664
$ws.Range("H86").Value = 14596.357
$ws.Range("I86").Value = 5848.484
$ws.Range("J86").Value = 39249.453
$ws.Range("K86").Value = 5848.484
$ws.Range("L86").Value = 39249.453
$ws.Range("M86").Value = -4725.484
$ws.Range("N86").Value = -41495.453
$ws.Range("H89").Value = 14596.357
$ws.Range("I89").Value = 5848.484
$ws.Range("J89").Value = 39249.453
$ws.Range("K89").Value = 29242.42
$ws.Range("L89").Value = 196247.265
$ws.Range("M89").Value = -23626.42
$ws.Range("N89").Value = -207479.265

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 55559680
$ws.Range("I31").Value = 1750
$ws.Range("J31").Value = 58827790
$ws.Range("K31").Value = 1750
$ws.Range("L31").Value = 58827790
$ws.Range("M31").Value = -1455
$ws.Range("N31").Value = -58828380
$ws.Range("H34").Value = 55559680
$ws.Range("I34").Value = 1750
$ws.Range("J34").Value = 58827790
$ws.Range("K34").Value = 1750
$ws.Range("L34").Value = 58827790
$ws.Range("M34").Value = -1548
$ws.Range("N34").Value = -58828194
$ws.Range("H99").Value = 8925.647999999999
$ws.Range("I99").Value = 5187.4165
$ws.Range("J99").Value = 10720
$ws.Range("K99").Value = 5187.4165
$ws.Range("L99").Value = 10720
$ws.Range("M99").Value = -3689.4165
$ws.Range("N99").Value = -13716
$ws.Range("H107").Value = 642
$ws.Range("J107").Value = 686.6667
$ws.Range("L107").Value = 686.6667
$ws.Range("N107").Value = -4526.6667
$ws.Range("H126").Value = 8925.647999999999
$ws.Range("I126").Value = 5187.4165
$ws.Range("J126").Value = 10720
$ws.Range("K126").Value = 15562.2495
$ws.Range("L126").Value = 32160
$ws.Range("M126").Value = -13092.2495
$ws.Range("N126").Value = -37100
$ws.Range("H132").Value = 2568.6667
$ws.Range("I132").Value = 2853.9565
$ws.Range("K132").Value = 8561.869499999999
$ws.Range("M132").Value = -6031.869499999999
$ws.Range("H134").Value = 2155.9395
$ws.Range("I134").Value = 1982.6666
$ws.Range("K134").Value = 5947.9998
$ws.Range("M134").Value = -3412.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1693.1904
$ws.Range("I102").Value = 1539.9474
$ws.Range("K102").Value = 1539.9474
$ws.Range("M102").Value = 82.05259999999998
$ws.Range("H132").Value = 3085.8667
$ws.Range("I132").Value = 3202
$ws.Range("J132").Value = 1460
$ws.Range("K132").Value = 9606
$ws.Range("L132").Value = 4380
$ws.Range("M132").Value = -7076
$ws.Range("N132").Value = -9440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4315.4287
$ws.Range("I40").Value = 3601.6
$ws.Range("K40").Value = 3601.6
$ws.Range("M40").Value = -3465.6
$ws.Range("H46").Value = 2011.3914
$ws.Range("I46").Value = 825.8570999999999
$ws.Range("K46").Value = 825.8570999999999
$ws.Range("M46").Value = -637.8570999999999
$ws.Range("H55").Value = 1336.2858
$ws.Range("J55").Value = 1409.1666
$ws.Range("L55").Value = 1409.1666
$ws.Range("N55").Value = -1755.1666
$ws.Range("H122").Value = 5697.6665
$ws.Range("I122").Value = 4400
$ws.Range("K122").Value = 13200
$ws.Range("M122").Value = -10750
$ws.Range("H132").Value = 37043790
$ws.Range("I132").Value = 3148.889
$ws.Range("J132").Value = 222247000
$ws.Range("K132").Value = 9446.667000000001
$ws.Range("L132").Value = 666741000
$ws.Range("M132").Value = -6916.667000000001
$ws.Range("N132").Value = -666746060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 946
$ws.Range("I14").Value = 946
$ws.Range("K14").Value = 946
$ws.Range("M14").Value = -778
$ws.Range("H100").Value = 37409450
$ws.Range("I100").Value = 48097188
$ws.Range("J100").Value = 2366.3333
$ws.Range("K100").Value = 96194376
$ws.Range("L100").Value = 4732.6666
$ws.Range("M100").Value = -96193835
$ws.Range("N100").Value = -5814.6666
$ws.Range("H132").Value = 1509.8043
$ws.Range("I132").Value = 1367.9286
$ws.Range("K132").Value = 4103.7858
$ws.Range("M132").Value = -1573.7858
